$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-10 Sunday" "2023-12-11 Monday"

Replace-Text "28×96=2688" "97×49=4753"
Replace-Text "26×34=884" "79×60=4740"
Replace-Text "32×69=2208" "69×70=4830"
Replace-Text "21×16=336" "28×29=812"
Replace-Text "36×22=792" "31×49=1519"

Replace-Text "13×40=520" "21×93=1953"
Replace-Text "46×67=3082" "58×62=3596"
Replace-Text "51×57=2907" "51×84=4284"
Replace-Text "98×18=1764" "23×23=529"
Replace-Text "73×43=3139" "82×55=4510"

Replace-Text "62×27=1674" "43×66=2838"
Replace-Text "26×63=1638" "17×47=799"
Replace-Text "64×66=4224" "85×16=1360"
Replace-Text "66×52=3432" "25×55=1375"
Replace-Text "99×61=6039" "11×51=561"

Replace-Text "37×90=3330" "38×20=760"
Replace-Text "57×53=3021" "70×31=2170"
Replace-Text "62×80=4960" "27×90=2430"
Replace-Text "78×67=5226" "84×19=1596"
Replace-Text "38×37=1406" "70×52=3640"

Replace-Text "12×94=1128" "14×25=350"
Replace-Text "93×44=4092" "73×30=2190"
Replace-Text "43×53=2279" "33×91=3003"
Replace-Text "69×93=6417" "19×83=1577"
Replace-Text "93×51=4743" "75×26=1950"

Write-Output "Done replacing all values"
